# Update the LR-pairs sheet (Wnt1-Fzd2) with refreshed TPM-derived values.
# The "Sending cluster" labels shift down one group (FAPs->ECs, MuSCs->FAPs)
# and every expression / specificity / weight column is recalculated
# accordingly, while "Target cluster" (column D) keeps the same cluster
# order (ECs, FAPs, MuSCs, Resolving-Mac) for each sender block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,7).Value = 0.131499
$ws.Cells.Item(2,8).Value = 0.394497
$ws.Cells.Item(2,9).Value = 0.3654391092296077
$ws.Cells.Item(2,10).Value = 0.3654391092296077
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.1503136666666667
$ws.Cells.Item(2,14).Value = 0.450941
$ws.Cells.Item(2,15).Value = 0.008029526741163598
$ws.Cells.Item(2,16).Value = 0.008029526741163598
$ws.Cells.Item(2,17).Value = 0.019766096853
$ws.Cells.Item(2,18).Value = 0.177894871677
$ws.Cells.Item(2,19).Value = 0.00293430309982614
$ws.Cells.Item(2,20).Value = 0.00293430309982614
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,7).Value = 0.131499
$ws.Cells.Item(3,8).Value = 0.394497
$ws.Cells.Item(3,9).Value = 0.3654391092296077
$ws.Cells.Item(3,10).Value = 0.3654391092296077
$ws.Cells.Item(3,15).Value = 0.8389317081486641
$ws.Cells.Item(3,16).Value = 0.8389317081486641
$ws.Cells.Item(3,17).Value = 2.065178425935
$ws.Cells.Item(3,18).Value = 18.586605833415
$ws.Cells.Item(3,19).Value = 0.306578456130321
$ws.Cells.Item(3,20).Value = 0.306578456130321
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,7).Value = 0.131499
$ws.Cells.Item(4,8).Value = 0.394497
$ws.Cells.Item(4,9).Value = 0.3654391092296077
$ws.Cells.Item(4,10).Value = 0.3654391092296077
$ws.Cells.Item(4,13).Value = 2.758466666666667
$ws.Cells.Item(4,14).Value = 8.2754
$ws.Cells.Item(4,15).Value = 0.1473530807662759
$ws.Cells.Item(4,16).Value = 0.1473530807662759
$ws.Cells.Item(4,17).Value = 0.3627356082
$ws.Cells.Item(4,18).Value = 3.2646204738
$ws.Cells.Item(4,19).Value = 0.05384857857746631
$ws.Cells.Item(4,20).Value = 0.05384857857746631
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,7).Value = 0.131499
$ws.Cells.Item(5,8).Value = 0.394497
$ws.Cells.Item(5,9).Value = 0.3654391092296077
$ws.Cells.Item(5,10).Value = 0.3654391092296077
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.1064366666666667
$ws.Cells.Item(5,14).Value = 0.31931
$ws.Cells.Item(5,15).Value = 0.005685684343896315
$ws.Cells.Item(5,16).Value = 0.005685684343896314
$ws.Cells.Item(5,17).Value = 0.01399631523
$ws.Cells.Item(5,18).Value = 0.12596683707
$ws.Cells.Item(5,19).Value = 0.002077771421994196
$ws.Cells.Item(5,20).Value = 0.002077771421994196
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,7).Value = 0.2283393333333333
$ws.Cells.Item(6,8).Value = 0.685018
$ws.Cells.Item(6,9).Value = 0.6345608907703922
$ws.Cells.Item(6,10).Value = 0.6345608907703922
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.1503136666666667
$ws.Cells.Item(6,14).Value = 0.450941
$ws.Cells.Item(6,15).Value = 0.008029526741163598
$ws.Cells.Item(6,16).Value = 0.008029526741163598
$ws.Cells.Item(6,17).Value = 0.03432252243755556
$ws.Cells.Item(6,18).Value = 0.308902701938
$ws.Cells.Item(6,19).Value = 0.005095223641337457
$ws.Cells.Item(6,20).Value = 0.005095223641337457
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,7).Value = 0.2283393333333333
$ws.Cells.Item(7,8).Value = 0.685018
$ws.Cells.Item(7,9).Value = 0.6345608907703922
$ws.Cells.Item(7,10).Value = 0.6345608907703922
$ws.Cells.Item(7,15).Value = 0.8389317081486641
$ws.Cells.Item(7,16).Value = 0.8389317081486641
$ws.Cells.Item(7,17).Value = 3.586046015501111
$ws.Cells.Item(7,18).Value = 32.27441413951
$ws.Cells.Item(7,19).Value = 0.532353252018343
$ws.Cells.Item(7,20).Value = 0.532353252018343
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,7).Value = 0.2283393333333333
$ws.Cells.Item(8,8).Value = 0.685018
$ws.Cells.Item(8,9).Value = 0.6345608907703922
$ws.Cells.Item(8,10).Value = 0.6345608907703922
$ws.Cells.Item(8,13).Value = 2.758466666666667
$ws.Cells.Item(8,14).Value = 8.2754
$ws.Cells.Item(8,15).Value = 0.1473530807662759
$ws.Cells.Item(8,16).Value = 0.1473530807662759
$ws.Cells.Item(8,17).Value = 0.6298664396888889
$ws.Cells.Item(8,18).Value = 5.6687979572
$ws.Cells.Item(8,19).Value = 0.09350450218880958
$ws.Cells.Item(8,20).Value = 0.09350450218880958
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,7).Value = 0.2283393333333333
$ws.Cells.Item(9,8).Value = 0.685018
$ws.Cells.Item(9,9).Value = 0.6345608907703922
$ws.Cells.Item(9,10).Value = 0.6345608907703922
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.1064366666666667
$ws.Cells.Item(9,14).Value = 0.31931
$ws.Cells.Item(9,15).Value = 0.005685684343896315
$ws.Cells.Item(9,16).Value = 0.005685684343896314
$ws.Cells.Item(9,17).Value = 0.02430367750888889
$ws.Cells.Item(9,18).Value = 0.21873309758
$ws.Cells.Item(9,19).Value = 0.003607912921902119
$ws.Cells.Item(9,20).Value = 0.003607912921902118
